$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Hello"
$ws.Range("B5").Value = "World"
$ws.Range("C5").Value = "2025-10-01T18:27:38.913Z"
